# Round volume calculation (column O, volumeEllipsoid) to one decimal place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$colIndex = 15  # Column O

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colIndex)
    $val = $cell.Value2
    if ($val -ne $null -and $val -is [double]) {
        $cell.Value2 = [Math]::Round($val, 1)
    }
}
